# Escalate sexting intensity in the IsabellaJourney sheet (column B = "Text")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IsabellaJourney")

$ws.Range("B4").Value  = "let go with me babe... I'm done holding back 🥵"
$ws.Range("B6").Value  = "I'm cumming... god I'm cumming so hard for you right now"
$ws.Range("B7").Value  = "I'm right there babe... every nerve in my body is on fire and I need you to watch me 🥵"
$ws.Range("B8").Value  = "oh my god 😏"
$ws.Range("B9").Value  = "you need to see this babe... I don't let anyone see me like this 🥵"
$ws.Range("B11").Value = "I'm so close to the edge and I need you to see what happens when I fall"
$ws.Range("B12").Value = "my body is arching off the bed and my legs are trembling... I can't stop 🥵"
$ws.Range("B13").Value = "I'm rubbing my pussy and going deeper with every stroke babe... I can hear how wet I am"
$ws.Range("B14").Value = "oh fuck 😏"
$ws.Range("B15").Value = "this is what you're making me do to myself babe... watch 🥵"
$ws.Range("B17").Value = "tell me what you want babe... I want to hear you say it while I'm touching myself like this"
$ws.Range("B18").Value = "I keep imagining you here, feeling your skin against mine, your breath on my neck... god I need it 🥵"
$ws.Range("B19").Value = "I'm touching myself right now and I can't believe how wet I already am... you did this babe"
$ws.Range("B20").Value = "wow... okay that hit deeper than I expected 😏"
$ws.Range("B21").Value = "see what you're doing to me babe... I couldn't keep this from you 🥵"
$ws.Range("B23").Value = "I want to show you what happens when I stop holding back... I think you can handle it 🥵"
$ws.Range("B24").Value = "my hand is drifting lower and I can feel myself getting wet just from this conversation... you're dangerous babe"
$ws.Range("B25").Value = "mm you liked that? good... because I can already feel my body reacting to the way you're looking at me 😏"
